# Applies the "MAJ mapping suite review de NRISS" update:
#  - Metadata!B5 (Title) becomes what used to be in Metadata!B4 (Name)
#  - Metadata!B4 (Name) is cleared
#  - Metadata!B8 (Date) is bumped to the new run timestamp
#  - The three "FREncounterDocument.participant.individual.extension.<x>"
#    entries in "Mapping Table 1" switch from dotted to colon notation
#    (FHIR extension slice syntax), i.e. ".executant" -> ":executant", etc.

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")

# Name's current value becomes the new Title value; Name's value is cleared.
$nameValue = $wsMeta.Range("B4").Value2
$wsMeta.Range("B5").Value = $nameValue
$wsMeta.Range("B4").Value = ""

# Refresh the metadata Date stamp.
$wsMeta.Range("B8").Value = "2026-01-07T15:20:53+00:00"

$wsMap1 = $wb.Worksheets.Item("Mapping Table 1")
$wsMap1.Range("D10").Value = "FREncounterDocument.participant.individual.extension:executant"
$wsMap1.Range("D11").Value = "FREncounterDocument.participant.individual.extension:author"
$wsMap1.Range("D12").Value = "FREncounterDocument.participant.individual.extension:informant"
